$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.3082
$ws.Range("H2").Value = 6.9246
$ws.Range("I2").Value = 0.6638288620319053
$ws.Range("J2").Value = 0.6638288620319053
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 16.6516102408
$ws.Range("R2").Value = 149.8644921672
$ws.Range("S2").Value = 0.3112576564754311
$ws.Range("T2").Value = 0.3112576564754311

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.3082
$ws.Range("H3").Value = 6.9246
$ws.Range("I3").Value = 0.6638288620319053
$ws.Range("J3").Value = 0.6638288620319053
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 16.4121114088
$ws.Range("R3").Value = 147.7090026792
$ws.Range("S3").Value = 0.3067808614929092
$ws.Range("T3").Value = 0.3067808614929092

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.3082
$ws.Range("H4").Value = 6.9246
$ws.Range("I4").Value = 0.6638288620319053
$ws.Range("J4").Value = 0.6638288620319053
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 2.449684196599999
$ws.Range("R4").Value = 22.0471577694
$ws.Range("S4").Value = 0.04579034406356504
$ws.Range("T4").Value = 0.04579034406356504

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.168901
$ws.Range("H5").Value = 3.506703
$ws.Range("I5").Value = 0.3361711379680947
$ws.Range("J5").Value = 0.3361711379680947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 8.432581172377333
$ws.Range("R5").Value = 75.893230551396
$ws.Range("S5").Value = 0.1576247231226878
$ws.Range("T5").Value = 0.1576247231226877

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.168901
$ws.Range("H6").Value = 3.506703
$ws.Range("I6").Value = 0.3361711379680947
$ws.Range("J6").Value = 0.3361711379680947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 8.311296004617333
$ws.Range("R6").Value = 74.801664041556
$ws.Range("S6").Value = 0.1553576188284911
$ws.Range("T6").Value = 0.1553576188284911

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.168901
$ws.Range("H7").Value = 3.506703
$ws.Range("I7").Value = 0.3361711379680947
$ws.Range("J7").Value = 0.3361711379680947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 1.240550345329666
$ws.Range("R7").Value = 11.164953107967
$ws.Range("S7").Value = 0.02318879601691588
$ws.Range("T7").Value = 0.02318879601691588

